{"js": "// Update the division-problem answers in the worksheet table.\n// The table has 20 rows; only every 4th row (0, 4, 8, 12, 16) holds the\n// 5 answer cells for that \"row\" of problems \u2014 the 3 rows in between are\n// blank work-space rows. We address each answer cell by its (row, col)\n// position so the edit is unambiguous regardless of any text collisions\n// between old/new values.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each entry: [tableRowIndex, columnIndex, newText]\nconst updates = [\n  [0, 0, \"39\u00f72=19, 1\"],\n  [0, 1, \"67\u00f79=7, 4\"],\n  [0, 2, \"68\u00f74=17, 0\"],\n  [0, 3, \"68\u00f79=7, 5\"],\n  [0, 4, \"54\u00f73=18, 0\"],\n\n  [4, 0, \"57\u00f72=28, 1\"],\n  [4, 1, \"46\u00f72=23, 0\"],\n  [4, 2, \"48\u00f78=6, 0\"],\n  [4, 3, \"31\u00f72=15, 1\"],\n  [4, 4, \"99\u00f75=19, 4\"],\n\n  [8, 0, \"19\u00f77=2, 5\"],\n  [8, 1, \"37\u00f73=12, 1\"],\n  [8, 2, \"42\u00f74=10, 2\"],\n  [8, 3, \"42\u00f74=10, 2\"],\n  [8, 4, \"95\u00f77=13, 4\"],\n\n  [12, 0, \"51\u00f75=10, 1\"],\n  [12, 1, \"40\u00f73=13, 1\"],\n  [12, 2, \"69\u00f72=34, 1\"],\n  [12, 3, \"68\u00f79=7, 5\"],\n  [12, 4, \"65\u00f77=9, 2\"],\n\n  [16, 0, \"34\u00f74=8, 2\"],\n  [16, 1, \"18\u00f73=6, 0\"],\n  [16, 2, \"52\u00f75=10, 2\"],\n  [16, 3, \"61\u00f78=7, 5\"],\n  [16, 4, \"34\u00f78=4, 2\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem answers in the worksheet table.\n# The table has 20 rows; only every 4th row (1, 5, 9, 13, 17 in Word's\n# 1-based Table.Cell indexing) holds the 5 answer cells for that \"row\" of\n# problems -- the 3 rows in between are blank work-space rows. Each answer\n# cell is addressed by its (row, col) position so the edit is unambiguous\n# regardless of any text collisions between old/new values.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each entry: row, column, new text (1-based, matching Table.Cell(r, c))\n$updates = @(\n    @(1, 1, \"39\u00f72=19, 1\"),\n    @(1, 2, \"67\u00f79=7, 4\"),\n    @(1, 3, \"68\u00f74=17, 0\"),\n    @(1, 4, \"68\u00f79=7, 5\"),\n    @(1, 5, \"54\u00f73=18, 0\"),\n\n    @(5, 1, \"57\u00f72=28, 1\"),\n    @(5, 2, \"46\u00f72=23, 0\"),\n    @(5, 3, \"48\u00f78=6, 0\"),\n    @(5, 4, \"31\u00f72=15, 1\"),\n    @(5, 5, \"99\u00f75=19, 4\"),\n\n    @(9, 1, \"19\u00f77=2, 5\"),\n    @(9, 2, \"37\u00f73=12, 1\"),\n    @(9, 3, \"42\u00f74=10, 2\"),\n    @(9, 4, \"42\u00f74=10, 2\"),\n    @(9, 5, \"95\u00f77=13, 4\"),\n\n    @(13, 1, \"51\u00f75=10, 1\"),\n    @(13, 2, \"40\u00f73=13, 1\"),\n    @(13, 3, \"69\u00f72=34, 1\"),\n    @(13, 4, \"68\u00f79=7, 5\"),\n    @(13, 5, \"65\u00f77=9, 2\"),\n\n    @(17, 1, \"34\u00f74=8, 2\"),\n    @(17, 2, \"18\u00f73=6, 0\"),\n    @(17, 3, \"52\u00f75=10, 2\"),\n    @(17, 4, \"61\u00f78=7, 5\"),\n    @(17, 5, \"34\u00f78=4, 2\")\n)\n\nforeach ($u in $updates) {\n    $t.Cell($u[0], $u[1]).Range.Text = $u[2]\n}\n"}
